$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as text (matches source
# workbook, where price/volume columns are inline strings, not numbers),
# then restore default General/Normal formatting so no stray number format is
# left applied to the cell.
function Set-TextValue($ws, $cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '66.900.19'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '3.513.82'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue $ws 'D5' '608.50'
$ws.Range('E5').Value = '  +0.66%  '
Set-TextValue $ws 'D6' '148.02'
$ws.Range('E6').Value = '  -1.80%  '
$ws.Range('D7').Value = '3.513.18'
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('E10').Value = '  -1.00%  '
Set-TextValue $ws 'D11' '7.99'
$ws.Range('E11').Value = '  +6.04%  '
$ws.Range('E12').Value = '  -1.85%  '
Set-TextValue $ws 'D13' '0.0000218'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '4.108.63'
$ws.Range('E14').Value = '  +0.93%  '
Set-TextValue $ws 'D15' '31.90'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '3.511.38'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').Value = '67.003.15'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('E18').Value = '  -0.38%  '
Set-TextValue $ws 'D19' '10.70'
$ws.Range('E19').Value = '  +8.26%  '
$ws.Range('E20').Value = '  -0.17%  '
Set-TextValue $ws 'D21' '15.35'
$ws.Range('E21').Value = '  -0.43%  '
Set-TextValue $ws 'D22' '438.32'
$ws.Range('E22').Value = '  -1.68%  '
Set-TextValue $ws 'D23' '0.609'
$ws.Range('E23').Value = '  -2.63%  '
Set-TextValue $ws 'D24' '79.41'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').Value = '3.656.15'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('E26').Value = '  +0.00%  '
Set-TextValue $ws 'D27' '0.0000122'
$ws.Range('E27').Value = '  -3.19%  '
Set-TextValue $ws 'D28' '9.77'
$ws.Range('E28').Value = '  -1.75%  '
Set-TextValue $ws 'D29' '8.28'
$ws.Range('E29').Value = '  -4.33%  '
Set-TextValue $ws 'D30' '2.52'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('E31').Value = '  -3.58%  '
$ws.Range('E32').Value = '  -1.11%  '
Set-TextValue $ws 'D33' '0.997'
$ws.Range('E33').Value = '  -0.19%  '
Set-TextValue $ws 'D34' '25.53'
$ws.Range('E34').Value = '  -0.14%  '
Set-TextValue $ws 'D35' '5.95'
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('E36').Value = '  -2.29%  '
Set-TextValue $ws 'D37' '8.03'
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('E38').Value = '  +0.01%  '
Set-TextValue $ws 'D39' '1.00'
$ws.Range('E39').Value = '  +0.11%  '
Set-TextValue $ws 'D40' '173.33'
$ws.Range('E40').Value = '  -2.35%  '
Set-TextValue $ws 'D41' '0.0893'
$ws.Range('E41').Value = '  -0.15%  '
Set-TextValue $ws 'D42' '5.41'
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').Value = '  -10.08%  '
Set-TextValue $ws 'D44' '0.894'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('E45').Value = '  -0.70%  '
Set-TextValue $ws 'D46' '27.80'
$ws.Range('E46').Value = '  -7.32%  '
Set-TextValue $ws 'D47' '1.27'
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D48' '2.46'
$ws.Range('E48').Value = '  -2.65%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws 'D49' '7.46'
$ws.Range('E49').Value = '  -1.69%  '
Set-TextValue $ws 'D50' '0.992'
$ws.Range('E50').Value = '  +0.74%  '
$ws.Range('E51').Value = '  -1.57%  '
